$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert a new row above the old row 4 ("Number of disability
#    persons" row). This shifts:
#      old row4 (Number of disability persons + values) -> row5
#      old row5 (Source citation, merged A5:H5)          -> row6
#    and creates a blank new row4 that we will populate with the
#    new "family with disabilities Persons" series.
# ---------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ---------------------------------------------------------------
# 2. Row 1 - new report title (merged A1:I1)
# ---------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Gurjaani Municipality"
$ws.Range("A1:I1").Merge()
$ws.Range("A1:I1").Font.Name = "Arial"
$ws.Range("A1:I1").Font.Size = 11
$ws.Range("A1:I1").Font.Bold = $true
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------
# 3. Row 2 - "(End of year, persons)" caption (text unchanged)
# ---------------------------------------------------------------
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.ThemeColor = 1
$ws.Range("A2").Interior.Pattern = 1
$ws.Range("A2").Interior.ThemeColor = 0
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------
# 4. Row 4 (new) - "family with disabilities Persons" series
# ---------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").Font.ThemeColor = 1
$ws.Range("A4").Interior.Pattern = 1
$ws.Range("A4").Interior.ThemeColor = 0
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true
$ws.Range("A4").Borders.Item(8).LineStyle = 1
$ws.Range("A4").Borders.Item(8).Weight = 2

$ws.Cells.Item(4,2).Value = 1400
$ws.Cells.Item(4,3).Value = 1358
$ws.Cells.Item(4,4).Value = 1307
$ws.Cells.Item(4,5).Value = 1351
$ws.Cells.Item(4,6).Value = 1336
$ws.Cells.Item(4,7).Value = 1323
$ws.Cells.Item(4,8).Value = 1305
$ws.Cells.Item(4,9).Value = 1328
$ws.Range("B4:I4").NumberFormat = "#\ ##0"
$ws.Range("B4:I4").Font.Name = "Arial"
$ws.Range("B4:I4").Font.Size = 10
$ws.Range("B4:I4").Font.ColorIndex = 1
$ws.Range("B4:I4").Interior.Pattern = 1
$ws.Range("B4:I4").Interior.ThemeColor = 0
$ws.Range("B4:I4").Borders.Item(8).LineStyle = 0
$ws.Range("B4:I4").Borders.Item(9).LineStyle = 0
$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------
# 5. Row 5 (old row4) - "disabilities Persons" series
# ---------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").Font.ThemeColor = 1
$ws.Range("A5").Interior.Pattern = 1
$ws.Range("A5").Interior.ThemeColor = 0
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Borders.Item(8).LineStyle = 0
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("A5").Borders.Item(9).Weight = 2

$ws.Cells.Item(5,2).Value = 1561
$ws.Cells.Item(5,3).Value = 1515
$ws.Cells.Item(5,4).Value = 1459
$ws.Cells.Item(5,5).Value = 1506
$ws.Cells.Item(5,6).Value = 1493
$ws.Cells.Item(5,7).Value = 1484
$ws.Cells.Item(5,8).Value = 1465
$ws.Cells.Item(5,9).Value = 1484
$ws.Range("B5:I5").NumberFormat = "#\ ##0"
$ws.Range("B5:I5").Font.Name = "Arial"
$ws.Range("B5:I5").Font.Size = 10
$ws.Range("B5:I5").Font.ColorIndex = 1
$ws.Range("B5:I5").Interior.Pattern = 1
$ws.Range("B5:I5").Interior.ThemeColor = 0
$ws.Range("B5:H5").Borders.Item(8).LineStyle = 0
$ws.Range("B5:H5").Borders.Item(9).LineStyle = 0
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2
$ws.Rows.Item(5).RowHeight = 21

# ---------------------------------------------------------------
# 6. Row 6 (old row5) - Source citation (merged A6:H6), text kept
# ---------------------------------------------------------------
$ws.Range("A6").Borders.Item(8).LineStyle = 0
$ws.Range("B6:H6").Borders.Item(8).LineStyle = 1
$ws.Range("B6:H6").Borders.Item(8).Weight = 2
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------
# 7. Column widths
# ---------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 20.81640625
